# Update countries & provincias Spain
#
# Source data refresh (Datos actualizados ... 16:52 -> 17:22): new case
# counts for several countries, and four country rows whose alphabetical/
# rank position moved in the source feed (so the row that used to show one
# country's numbers now shows another's). Since both the country label
# (col A) and its stats (cols B:H) live on the same row, each reordered
# pair below is expressed as two row rewrites: the row keeps its position
# but gets the label+numbers of whichever country now belongs there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados" banner ------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 17:22"

# --- Plain numeric refreshes (country stays put, only B:H change) ---------

# Estados Unidos (row 4)
$ws.Range("B4").Value = 290920
$ws.Range("C4").Value = 13759
$ws.Range("D4").Value = 14348
$ws.Range("E4").Value = 268728
$ws.Range("F4").Value = 6199
$ws.Range("G4").Value = 440
$ws.Range("H4").Value = 7844

# Alemania (row 7)
$ws.Range("B7").Value = 92150
$ws.Range("C7").Value = 991
$ws.Range("E7").Value = 64420
$ws.Range("G7").Value = 55
$ws.Range("H7").Value = 1330

# Suiza (row 13)
$ws.Range("E13").Value = 14791
$ws.Range("G13").Value = 50
$ws.Range("H13").Value = 641

# Brasil (row 20)
$ws.Range("B20").Value = 9244
$ws.Range("C20").Value = 50
$ws.Range("E20").Value = 8751
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 366

# Suecia (row 22)
$ws.Range("F22").Value = 379

# Grecia (row 45)
$ws.Range("B45").Value = 1673
$ws.Range("C45").Value = 60
$ws.Range("E45").Value = 1527
$ws.Range("G45").Value = 5
$ws.Range("H45").Value = 68

# Moldavia (row 68)
$ws.Range("D68").Value = 29
$ws.Range("E68").Value = 711
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 12

# --- Reordered rows (label + stats both change) ----------------------------

# rows 23/24 swap: Noruega now ranks above Australia
$ws.Range("A23").Value = "Noruega"
$ws.Range("C23").Value = 180
$ws.Range("D23").Value = 32
$ws.Range("E23").Value = 5456
$ws.Range("F23").Value = 98
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 62

$ws.Range("A24").Value = "Australia"
$ws.Range("B24").Value = 5550
$ws.Range("C24").Value = 96
$ws.Range("D24").Value = 585
$ws.Range("E24").Value = 4935
$ws.Range("F24").Value = 85
$ws.Range("H24").Value = 30

# rows 102/103 swap: Estado de Palestina now ranks above Malta
$ws.Range("A102").Value = "Estado de Palestina"
$ws.Range("B102").Value = 216
$ws.Range("C102").Value = 22
$ws.Range("D102").Value = 21
$ws.Range("E102").Value = 194
$ws.Range("F102").Value = 0
$ws.Range("H102").Value = 1

$ws.Range("A103").Value = "Malta"
$ws.Range("B103").Value = 213
$ws.Range("C103").Value = 11
$ws.Range("D103").Value = 2
$ws.Range("E103").Value = 211
$ws.Range("F103").Value = 3
$ws.Range("H103").Value = 0

# rows 119/120 swap: Isla de Man now ranks above Kenia
$ws.Range("A119").Value = "Isla de Man"
$ws.Range("C119").Value = 12
$ws.Range("D119").Value = 0
$ws.Range("E119").Value = 125
$ws.Range("F119").Value = 0
$ws.Range("H119").Value = 1

$ws.Range("A120").Value = "Kenia"
$ws.Range("B120").Value = 126
$ws.Range("C120").Value = 4
$ws.Range("D120").Value = 4
$ws.Range("E120").Value = 118
$ws.Range("F120").Value = 2
$ws.Range("H120").Value = 4

# rows 156/157/158 rotate: Birmania now ranks above Gabon and Haiti
$ws.Range("A156").Value = "Birmania"
$ws.Range("C156").Value = 1
$ws.Range("D156").Value = 0
$ws.Range("E156").Value = 20

$ws.Range("A157").Value = "Gabon"
$ws.Range("B157").Value = 21
$ws.Range("C157").Value = 0
$ws.Range("H157").Value = 1

$ws.Range("A158").Value = "Haiti"
$ws.Range("C158").Value = 2
$ws.Range("H158").Value = 0
